# Mark every feature as selected for the workshop ML input (column C)
# by flipping the flag from 0 to 1 for all data rows (rows 2-83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C83").Value = 1
